$d = $word.ActiveDocument

# --- Step 1: seed three fresh bulleted-list definitions (numId N, N+1, N+2) ---
# Word mints a brand-new numId each time a list template is applied to a
# paragraph, so we "burn" three of them on throw-away paragraphs appended
# at the very end of the body, remember the first minted numId, then
# delete the scratch paragraphs again (the numbering definitions stay
# behind in numbering.xml once minted).
$gallery = $word.ListGalleries.Item(1)
$bulletTemplate = $gallery.ListTemplates.Item(1)

$scratchStart = $d.Paragraphs.Count
for ($i = 0; $i -lt 3; $i++) {
    $end = $d.Range($d.Content.End - 1, $d.Content.End - 1)
    $end.InsertParagraphAfter()
    $cnt = $d.Paragraphs.Count
    $scratchPara = $d.Paragraphs.Item($cnt - 1)
    $scratchPara.Range.Text = "scratch"
    $scratchPara.Range.ListFormat.ApplyListTemplate($bulletTemplate)
}

# Remove the three scratch paragraphs again; the minted numId/abstractNum
# entries remain usable in numbering.xml.
$cnt = $d.Paragraphs.Count
$firstScratch = $d.Paragraphs.Item($cnt - 3)
$lastScratch = $d.Paragraphs.Item($cnt - 1)
$scratchRange = $d.Range($firstScratch.Range.Start, $lastScratch.Range.End)
$scratchRange.Delete()

# --- Step 2: insert the new content (paragraphs + table) right after the
# "Size: 204 cm" paragraph, i.e. just before the trailing empty paragraph ---
$target = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$insertPoint = $d.Range($target.Range.End, $target.Range.End)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/><w:p/><w:p/><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>1. Classic Cross-Country Ski Length</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="19"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Basic Formula</w:t></w:r><w:r><w:t xml:space="preserve">: For classic skiing, the ski length is usually around </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>110-120% of the skier's height</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="19"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Formula</w:t></w:r><w:r><w:t>: Classic Ski Length=Skier’s Height (cm)×1.1 to 1.2\</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>text{</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>Classic Ski Length} = \text{Skier's Height (cm)} \times 1.1 \text{ to } 1.2Classic Ski Length=Skier’s Height (cm)×1.1 to 1.2</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="19"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Skill Level Adjustment</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="19"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Beginner</w:t></w:r><w:r><w:t>: Aim for the shorter end (110-115% of height) for better control.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="19"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Intermediate to Advanced</w:t></w:r><w:r><w:t>: Use the longer end (115-120% of height) for more speed and glide.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="19"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Weight Adjustment</w:t></w:r><w:r><w:t xml:space="preserve">: If a skier is heavier or lighter than average, adjust by about </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>5 cm</w:t></w:r><w:r><w:t xml:space="preserve"> up or down.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Example</w:t></w:r><w:r><w:t>: For a skier who is 170 cm tall:</w:t></w:r></w:p><w:p><w:r><w:t>Classic Ski Length=170×1.1 to 1.2=187 to 204 cm\</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>text{</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>Classic Ski Length} = 170 \times 1.1 \text{ to } 1.2 = 187 \text{ to } 204 \text{ cm}Classic Ski Length=170×1.1 to 1.2=187 to 204 cm</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>2. Skate Cross-Country Ski Length</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="20"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Basic Formula</w:t></w:r><w:r><w:t xml:space="preserve">: Skate skis are generally shorter, around </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>105-110% of the skier's height</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="20"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Formula</w:t></w:r><w:r><w:t>: Skate Ski Length=Skier’s Height (cm)×1.05 to 1.1\</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>text{</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>Skate Ski Length} = \text{Skier's Height (cm)} \times 1.05 \text{ to } 1.1Skate Ski Length=Skier’s Height (cm)×1.05 to 1.1</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="20"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Skill Level Adjustment</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="20"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Beginner</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Opt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for 105% of height for better maneuverability.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="20"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Advanced</w:t></w:r><w:r><w:t>: Closer to 110% of height for better glide.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="20"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Weight Adjustment</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Similar to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> classic, add or subtract </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>5 cm</w:t></w:r><w:r><w:t xml:space="preserve"> for a heavier or lighter skier.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Example</w:t></w:r><w:r><w:t>: For a 170 cm tall skier:</w:t></w:r></w:p><w:p><w:r><w:t>Skate Ski Length=170×1.05 to 1.1=178.5 to 187 cm\</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>text{</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>Skate Ski Length} = 170 \times 1.05 \text{ to } 1.1 = 178.5 \text{ to } 187 \text{ cm}Skate Ski Length=170×1.05 to 1.1=178.5 to 187 cm</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>3. Adjustments Based on Terrain and Conditions</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Groomed Trails</w:t></w:r><w:r><w:t>: Stick closely to the recommended length from the formula.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Backcountry</w:t></w:r><w:r><w:t>: For more challenging or mixed terrain, shorter skis provide better control and stability.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>4. Quick Reference Table</w:t></w:r></w:p><w:p><w:r><w:t>Here’s a simplified table you could add to your app for easier reference:</w:t></w:r></w:p><w:tbl><w:tblPr><w:tblW w:w="0" w:type="auto"/><w:tblCellSpacing w:w="15" w:type="dxa"/><w:tblBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideH w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideV w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tblBorders><w:tblCellMar><w:top w:w="15" w:type="dxa"/><w:left w:w="15" w:type="dxa"/><w:bottom w:w="15" w:type="dxa"/><w:right w:w="15" w:type="dxa"/></w:tblCellMar><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="1245"/><w:gridCol w:w="2398"/><w:gridCol w:w="2232"/></w:tblGrid><w:tr><w:trPr><w:tblHeader/><w:tblCellSpacing w:w="15" w:type="dxa"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Height (cm)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Classic Ski Length (cm)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Skate Ski Length (cm)</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:tblCellSpacing w:w="15" w:type="dxa"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>150</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>165 - 180</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>158 - 165</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:tblCellSpacing w:w="15" w:type="dxa"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>160</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>176 - 192</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>168 - 176</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:tblCellSpacing w:w="15" w:type="dxa"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>170</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>187 - 204</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>178 - 187</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:tblCellSpacing w:w="15" w:type="dxa"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>180</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>198 - 216</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>189 - 198</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:tblCellSpacing w:w="15" w:type="dxa"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>190</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>209 - 228</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="0" w:type="auto"/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:r><w:t>199 - 209</w:t></w:r></w:p></w:tc></w:tr></w:tbl></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint.InsertXML($xml)
Write-Output "Paragraphs: $($d.Paragraphs.Count)  Tables: $($d.Tables.Count)"
